# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.744.01'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.659.40'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = "'303.00"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = "'0.3813"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').Value = "'0.3618"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').Value = "'51.23"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').Value = "'0.08191"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = "'1.233"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = "'1.001"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = "'6.466"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '1.663.83'
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('D18').Value = "'97.74"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('D19').Value = "'0.07012"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = "'6.812"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.60%  '
$ws.Range('D21').Value = "'17.63"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = "'12.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('D24').Value = '23.742.89'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = "'2.502"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('D26').Value = "'3.009"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').Value = "'153.24"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').Value = "'5.231"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').Value = "'134.21"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').Value = '1.845.20'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').Value = "'7.188"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.35%  '
$ws.Range('D33').Value = "'2.246"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('E34').Value = '  +4.92%  '
$ws.Range('D35').Value = "'1.055"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.32%  '
$ws.Range('D36').Value = "'0.02820"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').Value = "'0.2519"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'6.120"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = "'0.08805"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = "'0.07012"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('D42').Value = "'0.7008"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').Value = "'1.336"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').Value = "'16.14"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('D45').Value = "'0.6519"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = "'0.9999"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').Value = "'2.306"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').Value = "'3.965"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = "'0.07928"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = "'128.43"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = "'1.185"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
